# Update workbook "上海-漫展信息" (Shanghai Expo/Event Info) per the commit
# "Update gh-pages to output generated at 456a3b4".
#
# This commit re-scrapes event data from bilibili show listings:
#  - Sheet "展览" (sheet1, Exhibitions): several "want-to-go" counts (column F) bumped up.
#  - Sheet "演出" (sheet2, Performances): three cancelled/old events (rows 8-10) were
#    removed from the listing (causing all following rows to shift up by 3), the
#    "want-to-go" counts (column F) were refreshed for several rows, and the
#    running index in column A was kept sequential.
#  - Sheet "本地生活" (sheet3, Local life): several "want-to-go" counts bumped up.
#  - Sheet "全部类型" (sheet4, All types - a date-sorted merge of the above):
#    "want-to-go" counts refreshed to match the corresponding source rows.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# ---------------------------------------------------------------------------
# Sheet 1 (展览): refresh "想去人数" (want-to-go count) values in column F.
# ---------------------------------------------------------------------------
$sheet1Updates = @{
    2  = 1618
    5  = 9282
    6  = 284
    9  = 677
    10 = 607
    12 = 170
    13 = 305
    16 = 1559
    17 = 1346
    20 = 1426
    21 = 104
    22 = 265
    26 = 75
    27 = 335
    28 = 335
    29 = 1089
    32 = 247
    33 = 230
    34 = 64
    36 = 619
    38 = 144
    39 = 81
    40 = 167
    41 = 134
    42 = 549
    44 = 708
    45 = 252
    46 = 50
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# ---------------------------------------------------------------------------
# Sheet 2 (演出): remove 3 obsolete rows (old rows 8-10: a cancelled concert,
# and two other entries no longer listed), which shifts subsequent rows up by
# three positions, then refresh several "想去人数" counts and repair the
# sequential index in column A.
# ---------------------------------------------------------------------------
$ws2.Rows("8:10").Delete()

# Column A holds a simple running index (0,1,2,...) independent of row
# content; restore it after the row shift.
for ($r = 2; $r -le 39; $r++) {
    $ws2.Cells.Item($r, 1).Value = $r - 1
}

$sheet2Updates = @{
    7  = 3
    12 = 678
    20 = 948
    23 = 253
    24 = 645
    26 = 267
    28 = 165
    32 = 162
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Cells.Item($row, 6).Value = $sheet2Updates[$row]
}

# ---------------------------------------------------------------------------
# Sheet 3 (本地生活): refresh "想去人数" values in column F.
# ---------------------------------------------------------------------------
$sheet3Updates = @{
    5 = 340
    7 = 2265
    8 = 3375
    9 = 49
}
foreach ($row in $sheet3Updates.Keys) {
    $ws3.Cells.Item($row, 6).Value = $sheet3Updates[$row]
}

# ---------------------------------------------------------------------------
# Sheet 4 (全部类型): refresh "想去人数" values in column F (this sheet is a
# date-sorted merge of sheets 1-3, and the same refreshed counts propagate
# here for the matching rows).
# ---------------------------------------------------------------------------
$sheet4Updates = @{
    2  = 1618
    5  = 9282
    6  = 340
    8  = 3375
    11 = 677
    13 = 607
    14 = 170
    15 = 305
    16 = 1559
    17 = 678
    18 = 1346
    20 = 49
    21 = 1426
    22 = 104
    23 = 265
    27 = 335
    28 = 335
    29 = 1089
    34 = 948
    35 = 247
    37 = 230
    39 = 253
    41 = 619
    42 = 645
    43 = 144
    44 = 267
    45 = 167
    46 = 134
    47 = 165
    48 = 549
    49 = 708
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
